$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 286 (shifts old rows 286:342 down to 287:343)
$ws.Rows.Item(286).Insert()

# Populate the new row 286 with the new data record
$ws.Cells.Item(286, 1).Value = 8
$ws.Cells.Item(286, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(286, 3).Value = "Coquimbo"
$ws.Cells.Item(286, 4).Value = 45275
$ws.Cells.Item(286, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(286, 5).Value = 4
$ws.Cells.Item(286, 6).Value = 100112001
$ws.Cells.Item(286, 7).Value = "Berenjena"
$ws.Cells.Item(286, 8).Value = "Sin especificar"
$ws.Cells.Item(286, 9).Value = "Primera"
$ws.Cells.Item(286, 10).Value = 500
$ws.Cells.Item(286, 11).Value = 11000
$ws.Cells.Item(286, 12).Value = 12000
$ws.Cells.Item(286, 13).Value = 11500
$ws.Cells.Item(286, 14).Value = "$/caja 50 unidades"
$ws.Cells.Item(286, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(286, 16).Value = 230
$ws.Cells.Item(286, 17).Value = 50
$ws.Cells.Item(286, 18).Value = "Hortaliza"

Write-Host "Row inserted and populated"
